$wb = $excel.ActiveWorkbook

# 1. Update the data validation list on DataReference!D2:D1048576:
#    GROUP -> EXPERIMENTERGROUP
$wsRef = $wb.Worksheets.Item("DataReference")
$dv = $wsRef.Range("D2:D1048576").Validation
$dv.Formula1 = '"IMAGE,DATASET,PROJECT,EXPERIMENTERGROUP,EXPERIMENTER,TABLE,KEY_VALUES,TAG,COMMENT,ROI,FILE"'

# 2. Rename header on PSFBeadsInput: fitting_rss_threshold -> fitting_r2_threshold
$wsInput = $wb.Worksheets.Item("PSFBeadsInput")
$wsInput.Range("I1").Value = "fitting_r2_threshold"

# 3. Rename headers on PSFBeadsKeyValues: fit_rss_* -> fit_r2_*
$wsKeyValues = $wb.Worksheets.Item("PSFBeadsKeyValues")
$wsKeyValues.Range("J1").Value = "fit_r2_z_mean"
$wsKeyValues.Range("K1").Value = "fit_r2_z_median"
$wsKeyValues.Range("L1").Value = "fit_r2_z_std"
$wsKeyValues.Range("M1").Value = "fit_r2_y_mean"
$wsKeyValues.Range("N1").Value = "fit_r2_y_median"
$wsKeyValues.Range("O1").Value = "fit_r2_y_std"
$wsKeyValues.Range("P1").Value = "fit_r2_x_mean"
$wsKeyValues.Range("Q1").Value = "fit_r2_x_median"
$wsKeyValues.Range("R1").Value = "fit_r2_x_std"

Write-Host "RSS -> R2 rename applied"
